$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row for the ActivityLED connector part: a Digikey link to the
# Stewart Connector SS-52100-001 (added as a new shared string in A19).
$ws.Range("A19").Value = "https://www.digikey.com/en/products/detail/stewart-connector/SS-52100-001/7902377"

# Turn the existing Digikey URL text in A17 (PowerLED part) into a real
# hyperlink, and apply the built-in Hyperlink style so it renders/serializes
# the same way as the other hyperlinked cells (e.g. A22).
$ws.Hyperlinks.Add($ws.Range("A17"), "https://www.digikey.com/en/products/detail/goford-semiconductor/G6N02L/13664832")
$ws.Range("A17").Style = "Hyperlink"

# Mirror the saved selection state.
$ws.Range("F23").Select()
